# Auto-generated edit script: updates column F (想去人数 / interest count) values
# across sheets per the commit diff (gh-pages data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 288
$ws.Range("F5").Value = 19
$ws.Range("F6").Value = 28
$ws.Range("F8").Value = 1106
$ws.Range("F9").Value = 362
$ws.Range("F11").Value = 291
$ws.Range("F12").Value = 7881
$ws.Range("F14").Value = 9236
$ws.Range("F15").Value = 69
$ws.Range("F18").Value = 454
$ws.Range("F24").Value = 259
$ws.Range("F26").Value = 36
$ws.Range("F29").Value = 1600
$ws.Range("F31").Value = 56
$ws.Range("F32").Value = 291
$ws.Range("F34").Value = 38
$ws.Range("F35").Value = 323
$ws.Range("F37").Value = 886
$ws.Range("F40").Value = 1393
$ws.Range("F41").Value = 398
$ws.Range("F42").Value = 296
$ws.Range("F43").Value = 257
$ws.Range("F45").Value = 259
$ws.Range("F46").Value = 36
$ws.Range("F48").Value = 75

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 29
$ws.Range("F5").Value = 91
$ws.Range("F7").Value = 21
$ws.Range("F19").Value = 13
$ws.Range("F20").Value = 341

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2711
$ws.Range("F4").Value = 328
$ws.Range("F5").Value = 183

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 288
$ws.Range("F3").Value = 328
$ws.Range("F4").Value = 183
$ws.Range("F6").Value = 19
$ws.Range("F8").Value = 1106
$ws.Range("F9").Value = 362
$ws.Range("F12").Value = 91
$ws.Range("F13").Value = 291
$ws.Range("F14").Value = 7881
$ws.Range("F15").Value = 9236
$ws.Range("F16").Value = 69
$ws.Range("F22").Value = 1600
$ws.Range("F24").Value = 56
$ws.Range("F25").Value = 291
$ws.Range("F27").Value = 38
$ws.Range("F28").Value = 323
$ws.Range("F31").Value = 886
$ws.Range("F35").Value = 1393
$ws.Range("F36").Value = 398
$ws.Range("F38").Value = 296
$ws.Range("F39").Value = 257
$ws.Range("F42").Value = 259
$ws.Range("F43").Value = 36
$ws.Range("F46").Value = 13
$ws.Range("F47").Value = 341
$ws.Range("F48").Value = 75
